# Add set function and comment out
# Update the rule classification for specific rows in column A.
# Rows listed in $toStay flip from CHARTER to STAY.
# Rows listed in $toCharter flip from STAY to CHARTER.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RuleValue($Worksheet, $Row, $Value) {
    $Worksheet.Range("A$Row").Value = $Value
}

$toStay = @(4, 31, 67, 68, 71, 119, 122)
$toCharter = @(22, 23, 26, 27, 93, 133, 134, 135, 137, 138, 140, 143, 175)

foreach ($r in $toStay) {
    Set-RuleValue $ws $r "STAY"
}

foreach ($r in $toCharter) {
    Set-RuleValue $ws $r "CHARTER"
}
